$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that are removed in the updated dataset
$clearCells = @("W2", "Y2", "AA2", "AC2", "W3", "Y3", "Z3", "AA3", "AC3")
foreach ($cell in $clearCells) {
    $ws.Range($cell).ClearContents()
}

# Update numeric values for the refreshed capital structure data
$values = @{
    "D2" = 0.008835000000000001
    "E2" = -0.011
    "F2" = 0.0146
    "G2" = 0.08647998460354119
    "H2" = 0.08647998460354119
    "I2" = 0.06770592763664357
    "J2" = 0.0417163469679503
    "K2" = 1397.5
    "L2" = 0.03361961123941493
    "M2" = 736.073
    "N2" = 0.05083657936902591
    "O2" = 0.526706976744186
    "P2" = 707.3729999999999
    "Q2" = 0.04885442565887618
    "R2" = 0.5061703041144902
    "S2" = 28.7
    "T2" = 0.03899069793349301
    "U2" = 7590.400000000001
    "V2" = 0.5242278578927013
    "X2" = 0.3034124936540563
    "Z2" = 0.412388316458032
    "AB2" = 0.08060250072307074
    "AD2" = 111568.8
    "AE2" = 0
    "AF2" = 111568.8
    "AG2" = 103978.4
    "AH2" = 0.8851294744859102
    "AI2" = 0.8807601917689047
    "AJ2" = 0.8777689232265384
    "AK2" = 0.8731597903645276
    "AL2" = 357.5
    "AM2" = 230.7
    "AN2" = 31.69478139825573
    "AO2" = 7.872447552447553
    "AP2" = 29.53847902048238
    "AQ2" = 12.19939315127872
    "D3" = 0.0246
    "E3" = -0.204
    "G3" = 0.05485156149595168
    "H3" = 0.05485156149595168
    "I3" = 0.05097031229919034
    "J3" = 0.02548515614959517
    "K3" = 38.9
    "L3" = 0.004999357409073384
    "M3" = 29.173
    "N3" = 0.02352851036373901
    "O3" = 0.7499485861182519
    "P3" = 0.473
    "Q3" = 0.0003814823776110976
    "R3" = 0.01215938303341902
    "S3" = 28.7
    "T3" = 0.983786377815103
    "U3" = 691.6
    "V3" = 0.5577869182998628
    "X3" = 0.1114626786129166
    "AB3" = 0.07715980228153028
    "AD3" = 890.1
    "AF3" = 890.1
    "AG3" = 198.5
    "AH3" = 0.417887323943662
    "AI3" = 0.2406391089242748
    "AJ3" = 0.1380005561735261
    "AK3" = 0.06600605194027866
    "AL3" = 162.6
    "AM3" = 162.6
    "AN3" = 1.72
    "AO3" = 2.439114391143911
    "AP3" = 0.3835748792270531
    "AQ3" = 2.439114391143911
    "D4" = -0.00693
    "E4" = 0.182
    "F4" = 0.0146
    "G4" = 0.09376387367922574
    "H4" = 0.09376387367922574
    "I4" = 0.07156006748157576
    "J4" = 0.05240204154609434
    "K4" = 1358.6
    "L4" = 0.04021073193831947
    "M4" = 706.9
    "N4" = 0.05339406161957203
    "O4" = 0.5203150301781245
    "P4" = 706.9
    "Q4" = 0.05339406161957203
    "R4" = 0.5203150301781245
    "S4" = 0
    "T4" = 0
    "U4" = 6898.8
    "V4" = 0.5210849516213093
    "W4" = 0.1206282685324123
    "X4" = 0.4953623086951959
    "Y4" = -0.3747340401627837
    "Z4" = 0.3351944776791649
    "AA4" = 0.01756487494536499
    "AB4" = 0.08404519916461117
    "AC4" = -0.06648032421924618
    "AD4" = 110678.7
    "AE4" = 0
    "AF4" = 110678.7
    "AG4" = 103779.9
    "AH4" = 0.893160799883794
    "AI4" = 0.9000141492863556
    "AJ4" = 0.8868621559538947
    "AK4" = 0.8940716222875437
    "AL4" = 194.9
    "AM4" = 68.10000000000001
    "AN4" = 36.86095384000533
    "AO4" = 12.40533606977938
    "AP4" = 34.56334510091254
    "AQ4" = 35.50367107195301
}
foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
